# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1393
    "F3"  = 2702
    "F4"  = 577
    "F6"  = 6631
    "F7"  = 699
    "F9"  = 14
    "F10" = 11
    "F11" = 63
    "F12" = 8
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
